$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 632, shifting rows 632:673 down to 633:674
$ws.Rows(632).Insert()

# Populate the newly inserted row with the new data point.
# The date column stores plain text (e.g. "2026/01/15"), so force a text
# number format before assigning to avoid Excel auto-converting it to a
# date serial, then reset the style back to Normal so no style index is
# left behind on the cell.
$ws.Range("A632").NumberFormat = "@"
$ws.Range("A632").Value = "2026/01/15"
$ws.Range("A632").Style = "Normal"

$ws.Range("B632").Value = "木"
$ws.Range("C632").Value = 23
$ws.Range("D632").Value = 37
